$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to reflect the new "through" date
$wb.Sheets.Item(1).Name = "Through 2022-11-02"

# Update the row label for November to reflect the new "through" date
$ws.Range("A12").Value = "November (through 11-02)"

# Updated November figures (row 12), columns B (2015) through I (2022)
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 5
$ws.Range("D12").Value = 7
$ws.Range("E12").Value = 9
$ws.Range("F12").Value = 3
$ws.Range("G12").Value = 13
$ws.Range("H12").Value = 14
$ws.Range("I12").Value = 7

# Updated Total figures (row 13), columns B (2015) through I (2022)
$ws.Range("B13").Value = 260
$ws.Range("C13").Value = 491
$ws.Range("D13").Value = 717
$ws.Range("E13").Value = 624
$ws.Range("F13").Value = 485
$ws.Range("G13").Value = 1070
$ws.Range("H13").Value = 1455
$ws.Range("I13").Value = 1408
